$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '248.40'
Set-TextValue 'D3' '22.48'
Set-TextValue 'D4' '5.400'
Set-TextValue 'D5' '0.05691'
Set-TextValue 'D6' '3.404'
Set-TextValue 'D7' '6.317'
Set-TextValue 'D8' '0.8053'
Set-TextValue 'D9' '0.9174'
Set-TextValue 'D10' '0.1405'
Set-TextValue 'D11' '0.07441'
Set-TextValue 'D12' '0.03133'
Set-TextValue 'D13' '0.03045'
Set-TextValue 'D14' '0.09379'
Set-TextValue 'D15' '3.863'
Set-TextValue 'D16' '0.001576'
Set-TextValue 'D17' '0.04775'
Set-TextValue 'D19' '0.0005850'
Set-TextValue 'E19' '18OneONEWorstin24h'
Set-TextValue 'D20' '0.006459'
Set-TextValue 'D21' '0.004991'
Set-TextValue 'D22' '0.001007'
Set-TextValue 'D24' '3.698'
Set-TextValue 'D25' '2.199'
Set-TextValue 'D26' '0.3256'
Set-TextValue 'D40' '0.04019'
Set-TextValue 'B41' 'KickToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006839'
Set-TextValue 'E41' '40KickTokenKICK'
Set-TextValue 'B42' 'BKEXToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1068'
Set-TextValue 'E42' '41BKEXTokenBKK'
Set-TextValue 'B43' 'CEJI'
Set-TextValue 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.002720'
Set-TextValue 'E43' '42CEJICEJI'
Set-TextValue 'D44' '0.007517'
Set-TextValue 'D45' '0.00005799'
Set-TextValue 'D47' '0.4990'
Set-TextValue 'D48' '0.2104'
